# Adds the new Wall_Growth_RFU / Suspension_RFU columns to the "Data" sheet
# and documents them on the "DataDictionary" sheet.

$wb = $excel.ActiveWorkbook

$dataSheet = $wb.Worksheets.Item("Data")
$dictSheet = $wb.Worksheets.Item("DataDictionary")

# --- New headers on the Data sheet (columns R and S) ---
$dataSheet.Range("R1").Value = "Wall_Growth_RFU"
$dataSheet.Range("S1").Value = "Suspension_RFU"

# New headers pick up the same look (wrapped, size-10 font) as the rest of row 1
$dataSheet.Range("Q1").Copy() | Out-Null
$dataSheet.Range("R1:S1").PasteSpecial(-4122) | Out-Null

# --- New RFU values for rows 2-17 ---
$wallGrowth = @(249.58, 1019.53, 1855.4, 2102.13, 209.5, 1816.17, 0.57, 0.22, 0.24, 0.26, 339.73, 357.71, 1224.23, 1570.72, 1558.8, 769.38)
$suspension = @(11.58, 3.755, 6.105, 12.855, 7.905, 3.035, 0.18, 0.19, 0.225, 0.18, 8.21, 8.375, 3.495, 2.73, 9.15, 136.38)

for ($i = 0; $i -lt $wallGrowth.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 18).Value = $wallGrowth[$i]
    $dataSheet.Cells.Item($row, 19).Value = $suspension[$i]
}

# --- New rows on the DataDictionary sheet describing the new variables ---
$dictSheet.Range("A19").Value = "Wall_Growth_RFU"
$dictSheet.Range("C19").Value = "RFU"
$dictSheet.Range("B19").Value = "Chlorophyll a reading from growth on the tube wall"

$dictSheet.Range("A20").Value = "Suspension_RFU"
$dictSheet.Range("C20").Value = "RFU"
$dictSheet.Range("B20").Value = "Chlorophyll a reading after growth was suspended"

# Match the "Variable" column formatting used by the other rows (A13:A18)
$dictSheet.Range("A18").Copy() | Out-Null
$dictSheet.Range("A19:A20").PasteSpecial(-4122) | Out-Null

# Leave the new Data columns selected, as they were right after being entered
$dataSheet.Activate()
$dataSheet.Range("R1:S1").Select() | Out-Null
